# Actualización automática 2025-07-03 11:34:50
# Update the "PRESUPUESTO" (column G) values on the "VENTA MENSUAL" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$ws.Cells.Item(3, 7).Value  = 500
$ws.Cells.Item(5, 7).Value  = 4000
$ws.Cells.Item(6, 7).Value  = 4000
$ws.Cells.Item(7, 7).Value  = 1500
$ws.Cells.Item(10, 7).Value = 6000
$ws.Cells.Item(12, 7).Value = 5000
$ws.Cells.Item(13, 7).Value = 6500
$ws.Cells.Item(14, 7).Value = 500
$ws.Cells.Item(15, 7).Value = 500
$ws.Cells.Item(16, 7).Value = 8000
$ws.Cells.Item(18, 7).Value = 4500
$ws.Cells.Item(21, 7).Value = 5000

# Row 22 holds the column total; recompute it as the sum of the detail rows.
$ws.Cells.Item(22, 7).Value = 52500
